$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.932399999999998
$ws.Range("E3").Value = 15.95489999999999
$ws.Range("A12").Value = -21.585
$ws.Range("B14").Value = 6.716599999999995
$ws.Range("B26").Value = 4.137600000000003
$ws.Range("E30").Value = 15.72660000000001
$ws.Range("B31").Value = 4.996300000000002
$ws.Range("A32").Value = -21.24560000000001
$ws.Range("B35").Value = 9.464600000000006
$ws.Range("A36").Value = -19.96060000000001
$ws.Range("B37").Value = 8.825800000000005
$ws.Range("A38").Value = -19.1903
$ws.Range("E44").Value = 16.76099999999999
$ws.Range("B45").Value = 6.916099999999997
$ws.Range("A46").Value = -21.74730000000001
$ws.Range("A54").Value = -21.83859999999999
$ws.Range("A55").Value = -22.50450000000001
$ws.Range("B57").Value = 4.774899999999995
$ws.Range("E58").Value = 16.21980000000001
$ws.Range("A67").Value = -21.43209999999998
$ws.Range("A69").Value = -21.56469999999997
$ws.Range("A72").Value = -22.07000000000002
$ws.Range("E84").Value = 16.6548
$ws.Range("E89").Value = 17.48280000000002
$ws.Range("A91").Value = -21.37300000000001
$ws.Range("E91").Value = 17.95030000000002
$ws.Range("E92").Value = 18.00540000000002
$ws.Range("A99").Value = -20.11719999999999
$ws.Range("B100").Value = 5.238399999999997
$ws.Range("B102").Value = 8.244800000000003
$ws.Range("E102").Value = 16.76679999999999
